$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.095.09"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.589.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "529.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.72"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.602.71"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.43"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.054.48"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.034.31"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.44"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.587.08"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "346.92"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.33"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0735"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.85"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.77"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.93"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.97"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.83"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.825"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.74"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "268.58"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.77%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.594"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0518"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.40"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.956.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.61"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0221"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.39%  "
